$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Numero de orden"
$ws.Range("B2").Value = "Titulo"
$ws.Range("C2").Value = "responsable (via mail) "
$ws.Range("D2").Value = "Horas aprox"
$ws.Range("E2").Value = "Notas"

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Muestra desde excel"
$ws.Range("C3").Value = "alicemarcelaramirez@gmail.com"
$ws.Range("D3").Value = 5
$ws.Range("E3").Value = "Notas desde excel..."

$ws.Columns.Item(1).ColumnWidth = 21.29
$ws.Columns.Item(2).ColumnWidth = 25.43
$ws.Columns.Item(3).ColumnWidth = 27.57
$ws.Columns.Item(5).ColumnWidth = 29.0
